# Add new "Viraland cholera" dataset rows (r_courses) to the table of tables.
# Five releases of cholera surveillance data for a fictional country
# (Viraland), used for an intro R course.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# double-quote character, used to build formula strings without having to
# fight PowerShell's escaping rules inside single-quoted literals.
$q = '"'

# Columns: A name | B type | C extension | D type_version | E data_version
#          F language | G country | H scale | I subject | J context
#          K fictional | L year | M description | N usage | O license
#          P group_identifier (formula) | Q unique_identifier (formula)

$newRows = @(
    @("viraland_cholera_20210525_linelist",     "linelist", "zip", 1, 1, "en", "zzz", "national", "cholera", "surveillance", "yes", 2021, "Clean cholera surveillance data for Viraland, 25 May 2021",  "r_courses", "CC by-NC-SA 4.0"),
    @("viraland_cholera_20210617_linelist",     "linelist", "zip", 1, 2, "en", "zzz", "national", "cholera", "surveillance", "yes", 2021, "Clean cholera surveillance data for Viraland, 17 June 2021", "r_courses", "CC by-NC-SA 4.0"),
    @("viraland_cholera_20210721_linelist",     "linelist", "zip", 1, 3, "en", "zzz", "national", "cholera", "surveillance", "yes", 2021, "Clean cholera surveillance data for Viraland, 21 July 2021", "r_courses", "CC by-NC-SA 4.0"),
    @("viraland_cholera_20210721_linelist_raw", "linelist", "zip", 2, 3, "en", "zzz", "national", "cholera", "surveillance", "yes", 2021, "Raw cholera surveillance data for Viraland, 21 July 2021",   "r_courses", "CC by-NC-SA 4.0"),
    @("viraland_cholera_20210721_labs",         "linelist", "zip", 3, 3, "en", "zzz", "national", "cholera", "surveillance", "yes", 2021, "Raw cholera lab data for Viraland, 21 July 2021",             "r_courses", "CC by-NC-SA 4.0")
)

$startRow = 56
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]

    $formulaP = '=_xlfn.CONCAT(SUBSTITUTE(I{0},{1} {1},{1}{1}),{1}_{1},J{0},{1}_{1},G{0},{1}_{1},L{0})' -f $r, $q
    $formulaQ = '=_xlfn.CONCAT(P{0},{1}_{1},B{0},{1}_{1},D{0},{1}_{1},E{0})' -f $r, $q

    $ws.Cells.Item($r, 16).Formula = $formulaP
    $ws.Cells.Item($r, 17).Formula = $formulaQ
}
